$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Constants" sheet: insert a new blank row above row 19, pushing the
#    existing Get-Transaction-Data / Init-state rows (and everything below)
#    down by one row.
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Rows("19:19").Insert()

# Reflect the author's new viewport/selection on the Constants sheet
# (scrolled up a little, new selection at B37) without changing which
# sheet is active in the workbook.
$wsConstants.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$wsConstants.Range("B37").Select()

# ---------------------------------------------------------------------------
# 2. "Messages" sheet: remove the now-unused spacer column B (values shift
#    left: old C -> B, old D -> C, old E -> D) and add the two new Immi
#    error messages used by the Get Transaction Data / Init state.
# ---------------------------------------------------------------------------
$wsMessages = $wb.Worksheets.Item("Messages")
$wsMessages.Columns("B:B").Delete()

$wsMessages.Range("A21").Value = "ImmiErrorVevoInput"
$wsMessages.Range("B21").Value = "An exception occurred while entering Vevo details on Immi website"

$wsMessages.Range("A22").Value = "ImmiErrorDownloadSearchResult"
$wsMessages.Range("B22").Value = "An exception occurred while retriving Vevo serach result from Immi website"

# Restore Messages as the active sheet/tab with the new selection.
$wsMessages.Activate()
$wsMessages.Range("A24").Select()
